$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1810.682651439943
$ws.Range("C2").Value = 1984.399082413055
$ws.Range("D2").Value = 2108.131502855095
$ws.Range("E2").Value = 2136.069371440038
$ws.Range("B3").Value = 1768.334091166898
$ws.Range("C3").Value = 1961.485732241424
$ws.Range("D3").Value = 2100.219135539186
$ws.Range("E3").Value = 2132.944830111488
$ws.Range("B4").Value = 1771.385937107427
$ws.Range("C4").Value = 1962.522550761016
$ws.Range("D4").Value = 2096.593753118568
$ws.Range("E4").Value = 2127.424284502383
$ws.Range("B5").Value = 1696.724731045043
$ws.Range("C5").Value = 1923.028305113323
$ws.Range("D5").Value = 2070.329385050662
$ws.Range("E5").Value = 2120.356950798305
$ws.Range("B6").Value = 1252.015568647364
$ws.Range("C6").Value = 1625.470157090052
$ws.Range("D6").Value = 1894.731092481184
$ws.Range("E6").Value = 2056.101009838984
$ws.Range("B7").Value = 1349.580071416546
$ws.Range("C7").Value = 1681.95337263134
$ws.Range("D7").Value = 1916.80201782857
$ws.Range("E7").Value = 2070.712256698685
$ws.Range("B8").Value = 1347.754846155285
$ws.Range("C8").Value = 1665.712715699316
$ws.Range("D8").Value = 1896.532347435891
$ws.Range("E8").Value = 2049.817509468166
$ws.Range("B9").Value = 1364.57542459144
$ws.Range("C9").Value = 1687.562082513365
$ws.Range("D9").Value = 1920.750419483241
$ws.Range("E9").Value = 2073.482490100606
$ws.Range("B10").Value = 1657.608033392178
$ws.Range("C10").Value = 1909.025596153625
$ws.Range("D10").Value = 2063.638446140928
$ws.Range("E10").Value = 2114.98496527085
$ws.Range("B11").Value = 1635.170310488515
$ws.Range("C11").Value = 1890.486000554413
$ws.Range("D11").Value = 2045.065604917479
$ws.Range("E11").Value = 2103.257864896071
$ws.Range("B12").Value = 1247.307099901423
$ws.Range("C12").Value = 1669.150632374289
$ws.Range("D12").Value = 1950.445220564697
$ws.Range("E12").Value = 2043.878760777674
$ws.Range("B13").Value = 1529.951663943509
$ws.Range("C13").Value = 1835.465197322763
$ws.Range("D13").Value = 2019.552879454643
$ws.Range("E13").Value = 2086.80489728298
